# Added CodigoDistribuidor from master to output
#
# The matched-product rows (originally A2:J15) are re-joined against the
# distributor master on "Matched Name" (column C). Because a given matched
# product can correspond to more than one CODIGODISTRIBUIDOR in the master,
# some rows fan out into several output rows (one per matching distributor
# code), which is why the sheet grows from 14 data rows to 32 data rows
# (A2:J33) while the set of unique product/category combinations stays the
# same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row values (in output order) for columns A-D, H, I.
# E is always "NO", F is always "NO MATCH", G is always FALSE and
# J is always the size-mismatch message for every data row, before and after.
$rows = @(
    @{ A=14006577; B='POETT MUSICA PRMAVERAL*4.5LT'; C='mar in  primavera multiuso'; D=55.31914893617022; H=$false; I=$true },
    @{ A=14001197; B='Desinfectante Poett Solo Para Ti 4.5Lt.'; C='desinfectante pino sapolio prof'; D=73.01587301587303; H=$true; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=11001791.3; B='B LIMPIADOR POETT BEBE UND BOT 880 ML'; C='limpiatodo bebe sapolio'; D=62.74509803921568; H=$false; I=$false },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR.'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR.'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR.'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=14006766; B='POETT MUSICA PRMAVERAL*4.5LT'; C='mar in  primavera multiuso'; D=55.31914893617022; H=$false; I=$true },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR. BONIF'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR. BONIF'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=14005952; B='LEJIA CLOROX T.ACC. TRAD.BOT.3X4000 GR. BONIF'; C='lej a cloro patito'; D=58.33333333333333; H=$false; I=$false },
    @{ A=11000273; B='BONIF:POETT MAXIGALON FRESCURA LAVANDA*4.5L'; C='lavav crema patito manzana'; D=58.62068965517242; H=$false; I=$false },
    @{ A=14001197; B='Desinf. Poett Primavera 4.5Lt.'; C='mar in  primavera kekera'; D=62.22222222222222; H=$false; I=$true },
    @{ A=11002101; B='POET 880ML BEBE*** POR 5 PQT.POE 880'; C='pulso rpto a bebe sap'; D=61.90476190476191; H=$false; I=$false },
    @{ A=11001846; B=$null; C='vaini  negrit'; D=40; H=$false; I=$false },
    @{ A=11002670; B='BONIF CLOROX 860 (POR 10 PAQ CLRX 668)'; C='ac 6000 anticorrosivo x12'; D=45.71428571428572; H=$false; I=$false },
    @{ A=14006766; B='POETT FLORES DE PRIMAVERA*4.5LT'; C='deter patito floral   pe'; D=59.57446808510638; H=$false; I=$false },
    @{ A=14006766; B='POETT FLORES DE PRIMAVERA*4.5LT'; C='deter patito floral   pe'; D=59.57446808510638; H=$false; I=$false },
    @{ A=14006766; B='POETT BEBE*4.5LT'; C='det pat bebe'; D=72.72727272727273; H=$false; I=$false },
    @{ A=14006766; B='POETT BEBE*4.5LT'; C='det pat bebe'; D=72.72727272727273; H=$false; I=$false },
    @{ A=14006766; B='POETT BEBE*4.5LT'; C='det pat bebe'; D=72.72727272727273; H=$false; I=$false },
    @{ A=14006766; B='POETT BEBE*4.5LT'; C='det pat bebe'; D=72.72727272727273; H=$false; I=$false },
    @{ A=11000310; B='LEJ.CLOROX TRIP.ACCION TRAD. (CJAx3) 1x4000GR'; C='lej a cloro patito'; D=62.5; H=$true; I=$false },
    @{ A=11000310; B='LEJ.CLOROX TRIP.ACCION TRAD. (CJAx3) 1x4000GR'; C='lej a cloro patito'; D=62.5; H=$true; I=$false },
    @{ A=11000310; B='LEJ.CLOROX TRIP.ACCION TRAD. (CJAx3) 1x4000GR'; C='lej a cloro patito'; D=62.5; H=$true; I=$false }
)

# Clear the old data block (A2:J15) before rewriting, since it is being
# fully replaced / expanded.
$ws.Range("A2:J15").Clear()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    if ($row.B -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row.B
    }
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = "NO"
    $ws.Cells.Item($r, 6).Value = "NO MATCH"
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = "❌ Size mismatch. No homologation allowed."
    $r++
}
